# Apply strikethrough formatting to the "Creare un nuovo produttore ..."
# bullet list block (the producer/consumer API requirements that are no
# longer relevant) and to the already-highlighted "Scalare il credito ..."
# bullet.
#
# These correspond to the bullets, in order, starting at "Creare un nuovo
# produttore specificando quanto necessario" through "Aggiornare i costi
# dell'energia associati a una o più fasce orarie del produttore.", plus
# the separate "Scalare il credito di un utente ..." bullet further down.

$d = $word.ActiveDocument

$startMarker = "Creare un nuovo produttore specificando quanto necessario"
$endMarker = "Aggiornare i costi dell"
$scalareMarker = "Scalare il credito di u"

$startIndex = -1
$endIndex = -1
$scalareIndex = -1

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($startIndex -eq -1 -and $t.StartsWith($startMarker)) {
        $startIndex = $i
    }
    if ($endIndex -eq -1 -and $t.StartsWith($endMarker)) {
        $endIndex = $i
    }
    if ($scalareIndex -eq -1 -and $t.StartsWith($scalareMarker)) {
        $scalareIndex = $i
    }
}

# Strike through each paragraph individually (rather than one big range)
# so that every paragraph mark picks up the formatting too, matching a
# per-paragraph application of the strikethrough.
for ($i = $startIndex; $i -le $endIndex; $i++) {
    $d.Paragraphs.Item($i).Range.Font.StrikeThrough = 1
}

$d.Paragraphs.Item($scalareIndex).Range.Font.StrikeThrough = 1

# The hyperlink run inside the "Emissione di CO2 ..." bullet sometimes
# needs to be targeted explicitly so its rPr also receives <w:strike/>.
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $h = $d.Hyperlinks.Item($i)
    if ($h.Address -eq "https://www.isprambiente.gov.it/files2020/pubblicazioni/rapporti/Rapporto317_2020.pdf") {
        $h.Range.Font.StrikeThrough = 1
    }
}
